# Modified TOF sensor pins
# Re-pair each TOF sensor's Chip-Enable (CE) pin with its GPIO1 pin on
# consecutive rows, and drop the redundant "TOF Sensor N Chip Enable"
# descriptions for sensors 2-4 (keep it only for sensor 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 27 (TOF1_CE / description) is unchanged.

# Row 28: was TOF2_CE / "TOF Sensor 2 Chip Enable" -> becomes TOF1_GPIO1, no description
$ws.Cells.Item(28, 3).Value2 = "TOF1_GPIO1"
$ws.Cells.Item(28, 4).ClearContents()

# Row 29: was TOF3_CE / "TOF Sensor 3 Chip Enable" -> becomes TOF2_CE, no description
$ws.Cells.Item(29, 3).Value2 = "TOF2_CE"
$ws.Cells.Item(29, 4).ClearContents()

# Row 30: was TOF4_CE / "TOF Sensor 4 Chip Enable" -> becomes TOF2_GPIO1, no description
$ws.Cells.Item(30, 3).Value2 = "TOF2_GPIO1"
$ws.Cells.Item(30, 4).ClearContents()

# Row 34: was TOF1_GPIO1 -> becomes TOF3_CE
$ws.Cells.Item(34, 3).Value2 = "TOF3_CE"

# Row 35: was TOF2_GPIO1 -> becomes TOF3_GPIO1
$ws.Cells.Item(35, 3).Value2 = "TOF3_GPIO1"

# Row 36: was TOF3_GPIO1 -> becomes TOF4_CE
$ws.Cells.Item(36, 3).Value2 = "TOF4_CE"

# Row 37 (TOF4_GPIO1) is unchanged.

# Reflect where the edits were made in the saved view/selection.
$ws.Activate()
$ws.Range("D29").Select()
